# The "蟻に学べ" post (row 157) was removed from the posts sheet.
# Deleting the entire row shifts every subsequent row up by one,
# turning the old rows 158:215 into the new rows 157:214 and shrinking
# the used range from A1:C215 down to A1:C214.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(157).Delete()
